$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header "time (JST)" -> "time (UTC)"
$ws.Range("A1").Value = "time (UTC)"

# Shift logged timestamps from JST (UTC+9) to UTC (-0.375 days)
$ws.Range("A2").Value = 44198.62783564815
$ws.Range("A3").Value = 44198.62783564815
$ws.Range("A4").Value = 44198.627847222226
$ws.Range("A5").Value = 44198.627847222226
$ws.Range("A6").Value = 44198.627858796295
$ws.Range("A7").Value = 44198.62787037037
$ws.Range("A8").Value = 44198.62788194444
$ws.Range("A9").Value = 44199.62783564815

# Add new value in row 5, column F, matching the row's default style
$ws.Range("D5").Copy()
$ws.Range("F5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F5").Value = 1.234
